# Auto-generated script applying market-data refresh to Sheets/Spriggan_Profits.xlsx
# Updates columns H-N (currentAveragePrice*, LevePrice*, LeveProfit*) for the rows
# touched by the latest scheduled-runner data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 17233
$ws.Range("I12").Value = 17233
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 17233
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -17063
$ws.Range("N12").Value = $null

$ws.Range("H17").Value = 2111510
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2111510
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6334530
$ws.Range("N17").Value = -6334866

$ws.Range("H33").Value = 804.0714
$ws.Range("I33").Value = 387.25
$ws.Range("J33").Value = 1359.8334
$ws.Range("K33").Value = 387.25
$ws.Range("L33").Value = 1359.8334
$ws.Range("M33").Value = -158.25
$ws.Range("N33").Value = -1817.8334

$ws.Range("H51").Value = 4788.5
$ws.Range("I51").Value = 2770.5715
$ws.Range("J51").Value = 5875.077
$ws.Range("K51").Value = 2770.5715
$ws.Range("L51").Value = 5875.077
$ws.Range("M51").Value = -2286.5715
$ws.Range("N51").Value = -6843.077

$ws.Range("H70").Value = 11253.267
$ws.Range("I70").Value = 2257
$ws.Range("J70").Value = 19125
$ws.Range("K70").Value = 6771
$ws.Range("L70").Value = 57375
$ws.Range("M70").Value = -6501
$ws.Range("N70").Value = -57915

$ws.Range("H73").Value = 11253.267
$ws.Range("I73").Value = 2257
$ws.Range("J73").Value = 19125
$ws.Range("K73").Value = 6771
$ws.Range("L73").Value = 57375
$ws.Range("M73").Value = -5835
$ws.Range("N73").Value = -59247

$ws.Range("H80").Value = 740.1111
$ws.Range("I80").Value = 677.5
$ws.Range("J80").Value = 790.2
$ws.Range("K80").Value = 2032.5
$ws.Range("L80").Value = 2370.6
$ws.Range("M80").Value = -1034.5
$ws.Range("N80").Value = -4366.6

$ws.Range("H83").Value = 740.1111
$ws.Range("I83").Value = 677.5
$ws.Range("J83").Value = 790.2
$ws.Range("K83").Value = 6097.5
$ws.Range("L83").Value = 7111.8
$ws.Range("M83").Value = -1105.5
$ws.Range("N83").Value = -17095.8

$ws.Range("H87").Value = 64750
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 64750
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 64750
$ws.Range("N87").Value = -67246

$ws.Range("H90").Value = 64750
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 64750
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 194250
$ws.Range("N90").Value = -206730

$ws.Range("H111").Value = 1067.4445
$ws.Range("I111").Value = 950.875
$ws.Range("J111").Value = 2000
$ws.Range("K111").Value = 2852.625
$ws.Range("L111").Value = 6000
$ws.Range("M111").Value = 214.375
$ws.Range("N111").Value = -12134

$ws.Range("H135").Value = 100000710
$ws.Range("I135").Value = 100000710
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 900006390
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -900003855
$ws.Range("N135").Value = $null

$ws.Range("H138").Value = 2642.034
$ws.Range("I138").Value = 1525.0385
$ws.Range("J138").Value = 3522.0908
$ws.Range("K138").Value = 4575.1155
$ws.Range("L138").Value = 10566.2724
$ws.Range("M138").Value = 564.8845000000001
$ws.Range("N138").Value = -20846.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4389.6665
$ws.Range("I45").Value = 4750.5
$ws.Range("J45").Value = 3668
$ws.Range("K45").Value = 4750.5
$ws.Range("L45").Value = 3668
$ws.Range("M45").Value = -4373.5
$ws.Range("N45").Value = -4422

$ws.Range("H74").Value = 43484156
$ws.Range("I74").Value = 45460116
$ws.Range("J74").Value = 13000
$ws.Range("K74").Value = 45460116
$ws.Range("L74").Value = 13000
$ws.Range("M74").Value = -45459242
$ws.Range("N74").Value = -14748

$ws.Range("H77").Value = 43484156
$ws.Range("I77").Value = 45460116
$ws.Range("J77").Value = 13000
$ws.Range("K77").Value = 227300580
$ws.Range("L77").Value = 65000
$ws.Range("M77").Value = -227296212
$ws.Range("N77").Value = -73736

$ws.Range("H122").Value = 3452.75
$ws.Range("I122").Value = 3412
$ws.Range("J122").Value = 3575
$ws.Range("K122").Value = 10236
$ws.Range("L122").Value = 10725
$ws.Range("M122").Value = -7786
$ws.Range("N122").Value = -15625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 17999.334
$ws.Range("I94").Value = 17999.334
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 17999.334
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -17548.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1555183.6
$ws.Range("I16").Value = 2174457.2
$ws.Range("J16").Value = 6999.5
$ws.Range("K16").Value = 2174457.2
$ws.Range("L16").Value = 6999.5
$ws.Range("M16").Value = -2174170.2
$ws.Range("N16").Value = -7573.5

$ws.Range("H31").Value = 10826.846
$ws.Range("I31").Value = 7769.077
$ws.Range("J31").Value = 13884.615
$ws.Range("K31").Value = 7769.077
$ws.Range("L31").Value = 13884.615
$ws.Range("M31").Value = -7474.077
$ws.Range("N31").Value = -14474.615

$ws.Range("H34").Value = 10826.846
$ws.Range("I34").Value = 7769.077
$ws.Range("J34").Value = 13884.615
$ws.Range("K34").Value = 7769.077
$ws.Range("L34").Value = 13884.615
$ws.Range("M34").Value = -7567.077
$ws.Range("N34").Value = -14288.615

$ws.Range("H86").Value = 5449.1665
$ws.Range("I86").Value = 4999.4
$ws.Range("J86").Value = 7698
$ws.Range("K86").Value = 4999.4
$ws.Range("L86").Value = 7698
$ws.Range("M86").Value = -3876.4
$ws.Range("N86").Value = -9944

$ws.Range("H89").Value = 5449.1665
$ws.Range("I89").Value = 4999.4
$ws.Range("J89").Value = 7698
$ws.Range("K89").Value = 24997
$ws.Range("L89").Value = 38490
$ws.Range("M89").Value = -19381
$ws.Range("N89").Value = -49722

$ws.Range("H105").Value = 1251139
$ws.Range("I105").Value = 1429530.2
$ws.Range("J105").Value = 2400
$ws.Range("K105").Value = 1429530.2
$ws.Range("L105").Value = 2400
$ws.Range("M105").Value = -1427783.2
$ws.Range("N105").Value = -5894

$ws.Range("H113").Value = 1555183.6
$ws.Range("I113").Value = 2174457.2
$ws.Range("J113").Value = 6999.5
$ws.Range("K113").Value = 2174457.2
$ws.Range("L113").Value = 6999.5
$ws.Range("M113").Value = -2172287.2
$ws.Range("N113").Value = -11339.5

$ws.Range("H122").Value = 2717.7856
$ws.Range("I122").Value = 2717.7856
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8153.3568
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5703.3568
$ws.Range("N122").Value = $null

$ws.Range("H132").Value = 20001874
$ws.Range("I132").Value = 22224138
$ws.Range("J132").Value = 1507.4
$ws.Range("K132").Value = 66672414
$ws.Range("L132").Value = 4522.200000000001
$ws.Range("M132").Value = -66669884
$ws.Range("N132").Value = -9582.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 111829.336
$ws.Range("I5").Value = 167082.33
$ws.Range("J5").Value = 1323.3334
$ws.Range("K5").Value = 501246.99
$ws.Range("L5").Value = 3970.0002
$ws.Range("M5").Value = -501134.99
$ws.Range("N5").Value = -4194.0002

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = $null
$ws.Range("N68").Value = $null

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = $null
$ws.Range("N71").Value = $null

$ws.Range("H92").Value = 653.8
$ws.Range("I92").Value = 345
$ws.Range("J92").Value = 859.6667
$ws.Range("K92").Value = 1035
$ws.Range("L92").Value = 2579.0001
$ws.Range("M92").Value = 213
$ws.Range("N92").Value = -5075.0001

$ws.Range("H135").Value = 111829.336
$ws.Range("I135").Value = 167082.33
$ws.Range("J135").Value = 1323.3334
$ws.Range("K135").Value = 1503740.97
$ws.Range("L135").Value = 11910.0006
$ws.Range("M135").Value = -1501205.97
$ws.Range("N135").Value = -16980.0006

$ws.Range("H137").Value = 5557349
$ws.Range("I137").Value = 14287278
$ws.Range("J137").Value = 1939.909
$ws.Range("K137").Value = 42861834
$ws.Range("L137").Value = 5819.727000000001
$ws.Range("M137").Value = -42856734
$ws.Range("N137").Value = -16019.727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15437.4
$ws.Range("I102").Value = 8911.286
$ws.Range("J102").Value = 30665
$ws.Range("K102").Value = 8911.286
$ws.Range("L102").Value = 30665
$ws.Range("M102").Value = -7289.286
$ws.Range("N102").Value = -33909

$ws.Range("H125").Value = 148333
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 148333
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 148333
$ws.Range("N125").Value = -153253

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4308.0835
$ws.Range("I7").Value = 4325.3
$ws.Range("J7").Value = 4222
$ws.Range("K7").Value = 4325.3
$ws.Range("L7").Value = 4222
$ws.Range("M7").Value = -4213.3
$ws.Range("N7").Value = -4446

$ws.Range("H46").Value = 1047.875
$ws.Range("I46").Value = 1248.5
$ws.Range("J46").Value = 981
$ws.Range("K46").Value = 1248.5
$ws.Range("L46").Value = 981
$ws.Range("M46").Value = -1060.5
$ws.Range("N46").Value = -1357

$ws.Range("H126").Value = 4308.0835
$ws.Range("I126").Value = 4325.3
$ws.Range("J126").Value = 4222
$ws.Range("K126").Value = 12975.9
$ws.Range("L126").Value = 12666
$ws.Range("M126").Value = -10505.9
$ws.Range("N126").Value = -17606

$ws.Range("H132").Value = 20887458
$ws.Range("I132").Value = 24488168
$ws.Range("J132").Value = 3339.8
$ws.Range("K132").Value = 73464504
$ws.Range("L132").Value = 10019.4
$ws.Range("M132").Value = -73461974
$ws.Range("N132").Value = -15079.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7099.1113
$ws.Range("I62").Value = 6000
$ws.Range("J62").Value = 7236.5
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 7236.5
$ws.Range("M62").Value = -5376
$ws.Range("N62").Value = -8484.5

$ws.Range("H65").Value = 7099.1113
$ws.Range("I65").Value = 6000
$ws.Range("J65").Value = 7236.5
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 36182.5
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -42422.5

$ws.Range("H96").Value = 1673.4
$ws.Range("I96").Value = 1729.6666
$ws.Range("J96").Value = 1649.2858
$ws.Range("K96").Value = 1729.6666
$ws.Range("L96").Value = 1649.2858
$ws.Range("M96").Value = -356.6666
$ws.Range("N96").Value = -4395.2858

$ws.Range("H113").Value = 1282.25
$ws.Range("I113").Value = 1348.0714
$ws.Range("J113").Value = 821.5
$ws.Range("K113").Value = 4044.2142
$ws.Range("L113").Value = 2464.5
$ws.Range("M113").Value = -1874.2142
$ws.Range("N113").Value = -6804.5

$ws.Range("H132").Value = 15153412
$ws.Range("I132").Value = 16668541
$ws.Range("J132").Value = 2126.6667
$ws.Range("K132").Value = 50005623
$ws.Range("L132").Value = 6380.000100000001
$ws.Range("M132").Value = -50003093
$ws.Range("N132").Value = -11440.0001

